# Apply "edit daftar lalu lintas" changes:
#  - rename sheet "DaftarLaluLintas_ubah" -> "DaftarLaluLintas_Edit"
#  - fix the noSK value for the first record ("SK/001/DIV" -> "sk001")
#  - fill in the (previously blank) PengawalInternal (col B) values for rows 4-10
#  - update the sheet's active selection (C2) / scroll position

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DaftarLaluLintas_ubah")
$ws.Name = "DaftarLaluLintas_Edit"

# Correct the noSK text for the first data row.
$ws.Range("C2").Value = "sk001"

# Populate PengawalInternal for the remaining rows, alternating between the
# two guards already used in rows 2 and 3.
$ws.Range("B4").Value = "Wildan Cahyono"
$ws.Range("B5").Value = "EYONO BIN CAS"
$ws.Range("B6").Value = "Wildan Cahyono"
$ws.Range("B7").Value = "EYONO BIN CAS"
$ws.Range("B8").Value = "Wildan Cahyono"
$ws.Range("B9").Value = "EYONO BIN CAS"
$ws.Range("B10").Value = "Wildan Cahyono"

# Move the selection/scroll position as left by the editor.
$ws.Range("C2").Select()
